# Updated cryptos list on Wed Feb 22 09:44:06 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '24.074.43'
$ws.Range("E2").Value = '  -3.32%  '

# Row 3
Set-TextValue "D3" '1.638.49'
$ws.Range("E3").Value = '  -3.15%  '

# Row 4
Set-TextValue "D4" '0.9998'
$ws.Range("E4").Value = '  -0.56%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue "D5" '307.09'
$ws.Range("E5").Value = '  -2.81%  '

# Row 6
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue "D6" '1.0000'
$ws.Range("E6").Value = '  -0.51%  '

# Row 7
Set-TextValue "D7" '0.3872'
$ws.Range("E7").Value = '  -2.22%  '

# Row 8
Set-TextValue "D8" '0.3828'
$ws.Range("E8").Value = '  -4.19%  '

# Row 9
Set-TextValue "D9" '1.0000'
$ws.Range("E9").Value = '  -0.48%  '

# Row 10
Set-TextValue "D10" '49.34'
$ws.Range("E10").Value = '  -6.01%  '

# Row 11
Set-TextValue "D11" '1.341'
$ws.Range("E11").Value = '  -7.12%  '

# Row 12
Set-TextValue "D12" '0.08547'
$ws.Range("E12").Value = '  -2.09%  '

# Row 13
Set-TextValue "D13" '23.51'
$ws.Range("E13").Value = '  -7.92%  '

# Row 14
Set-TextValue "D14" '7.050'
$ws.Range("E14").Value = '  -4.75%  '

# Row 15
Set-TextValue "D15" '0.00001277'
$ws.Range("E15").Value = '  -4.73%  '

# Row 16
Set-TextValue "D16" '7.429'
$ws.Range("E16").Value = '  -5.63%  '

# Row 17
Set-TextValue "D17" '1.635.64'
$ws.Range("E17").Value = '  -9.83%  '

# Row 18
Set-TextValue "D18" '95.02'
$ws.Range("E18").Value = '  +0.11%  '

# Row 19
Set-TextValue "D19" '0.06875'
$ws.Range("E19").Value = '  -4.62%  '

# Row 20
Set-TextValue "D20" '20.75'
$ws.Range("E20").Value = '  +1.60%  '

# Row 21
Set-TextValue "D21" '6.868'
$ws.Range("E21").Value = '  -4.46%  '

# Row 22
Set-TextValue "D22" '0.9998'
$ws.Range("E22").Value = '  -0.72%  '

# Row 23
$ws.Range("E23").Value = '  -5.07%  '

# Row 24
Set-TextValue "D24" '24.075.66'
$ws.Range("E24").Value = '  -3.29%  '

# Row 25
Set-TextValue "D25" '2.327'
$ws.Range("E25").Value = '  -3.07%  '

# Row 26
$ws.Range("E26").Value = '  -6.50%  '

# Row 27
Set-TextValue "D27" '22.27'
$ws.Range("E27").Value = '  -3.72%  '

# Row 28
Set-TextValue "D28" '156.78'
$ws.Range("E28").Value = '  -3.35%  '

# Row 29
$ws.Range("E29").Value = '  +7.37%  '

# Row 30
Set-TextValue "D30" '139.66'
$ws.Range("E30").Value = '  -6.29%  '

# Row 31
Set-TextValue "D31" '5.343'
$ws.Range("E31").Value = '  -11.44%  '

# Row 32
Set-TextValue "D32" '2.411'
$ws.Range("E32").Value = '  -7.66%  '

# Row 33
Set-TextValue "D33" '1.817.45'
$ws.Range("E33").Value = '  -1.65%  '

# Row 34
Set-TextValue "D34" '6.837'
$ws.Range("E34").Value = '  -2.97%  '

# Row 35
Set-TextValue "D35" '0.07974'
$ws.Range("E35").Value = '  -6.35%  '

# Row 36
Set-TextValue "D36" '0.02869'
$ws.Range("E36").Value = '  -7.78%  '

# Row 37
Set-TextValue "D37" '0.2659'
$ws.Range("E37").Value = '  -7.38%  '

# Row 38
Set-TextValue "D38" '0.9430'
$ws.Range("E38").Value = '  -9.02%  '

# Row 39
Set-TextValue "D39" '0.09148'
$ws.Range("E39").Value = '  -5.59%  '

# Row 40
Set-TextValue "D40" '1.440'
$ws.Range("E40").Value = '  -2.26%  '

# Row 41
Set-TextValue "D41" '9.826'
$ws.Range("E41").Value = '  -8.80%  '

# Row 42
Set-TextValue "D42" '0.7494'
$ws.Range("E42").Value = '  -7.48%  '

# Row 43
Set-TextValue "D43" '12.95'
$ws.Range("E43").Value = '  -6.98%  '

# Row 44
Set-TextValue "D44" '15.91'
$ws.Range("E44").Value = '  -6.33%  '

# Row 45
$ws.Range("E45").Value = '  -6.11%  '

# Row 46
Set-TextValue "D46" '2.443'
$ws.Range("E46").Value = '  -7.21%  '

# Row 47
Set-TextValue "D47" '4.079'
$ws.Range("E47").Value = '  -3.29%  '

# Row 48
Set-TextValue "D48" '0.9994'
$ws.Range("E48").Value = '  -0.40%  '

# Row 49
Set-TextValue "D49" '0.08332'
$ws.Range("E49").Value = '  -6.85%  '

# Row 50
Set-TextValue "D50" '1.253'
$ws.Range("E50").Value = '  -9.84%  '

# Row 51
Set-TextValue "D51" '132.07'
$ws.Range("E51").Value = '  -5.33%  '

